$d = $word.ActiveDocument

$pairs = @(
    @("53÷4=", "27÷3="),
    @("31÷7=", "39÷3="),
    @("11÷4=", "57÷3="),
    @("85÷2=", "39÷9="),
    @("82÷9=", "83÷2="),
    @("89÷2=", "11÷7="),
    @("74÷4=", "41÷4="),
    @("94÷2=", "44÷4="),
    @("61÷8=", "91÷6="),
    @("96÷9=", "45÷3="),
    @("76÷5=", "30÷3="),
    @("42÷6=", "12÷9="),
    @("88÷2=", "79÷7="),
    @("84÷8=", "33÷5="),
    @("72÷4=", "88÷7="),
    @("40÷5=", "27÷6="),
    @("65÷8=", "60÷6="),
    @("97÷9=", "44÷2="),
    @("15÷3=", "46÷8="),
    @("73÷8=", "15÷8="),
    @("21÷3=", "80÷3="),
    @("96÷7=", "68÷9="),
    @("13÷9=", "51÷5="),
    @("59÷2=", "56÷4="),
    @("35÷4=", "31÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
